# TM59 test-fixture regeneration:
#  - report re-run with a new output dir override (/c/e), new "Author" (jovyan)
#    and a new run timestamp/date (2022-03-04)
#  - the small "readme" index table's columns were reordered to
#    index, sheet_name, Author, JobNo, Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")

# ---- Reorder / relabel the header row of the index table -------------
$ws.Range("B1").Value = "sheet_name"
$ws.Range("C1").Value = "Author"
$ws.Range("D1").Value = "JobNo"
$ws.Range("E1").Value = "Date"

# ---- Rewrite each data row in the new column order --------------------
$sheetNames = @(
    "Project Information",
    "Criterion Definitions",
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $sheetNames[$i]
    $ws.Cells.Item($row, 3).Value = "jovyan"
    $ws.Cells.Item($row, 4).Value = "/c/e"

    # "20220304" reads as a number to Excel's auto-typing; force it to stay
    # text (matching the original report's string column) via a formula
    # that yields a text result, then bake the formula down to a literal
    # value with Paste Special values-only so no formula is left behind.
    $ws.Cells.Item($row, 5).Formula = '="20220304"'
    $ws.Cells.Item($row, 5).Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

# ---- Bump the run timestamp on the "Project Information" sheet --------
$info = $wb.Worksheets.Item("Project Information")
$info.Range("B11").Value = "2022-03-04 17:31:00.947870"
